# Notas da prova 03
# Fill in column D ("Prova 03") grades for each student row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows without a grade yet (marked with "-")
$dashRows = @(5, 13, 24)
foreach ($r in $dashRows) {
    $ws.Range("D$r").Value = "-"
}

# Rows with new Prova 03 formulas (sum of six component grades)
$formulas = @{
    6  = "=12+7+15+7+6+20"
    7  = "=15+18+15+12+15+20"
    8  = "=15+20+15+15+15+20"
    9  = "=15+13+15+8+7+16"
    10 = "=12+15+8+8+15+20"
    11 = "=14+15+15+7+15+20"
    12 = "=15+15+15+15+12+20"
    14 = "=14+20+15+15+15+20"
    15 = "=3+12+12+7+5+12"
    16 = "=8+10+10+10+7+20"
    17 = "=10+15+12+12+7+18"
    18 = "=15+20+15+15+15+20"
    19 = "=10+15+15+15+15+20"
    20 = "=15+15+15+15+12+20"
    21 = "=15+15+15+12+14+20"
    22 = "=15+20+15+15+15+20"
    23 = "=10+20+15+15+12+20"
    25 = "=15+20+15+7+15+20"
    26 = "=15+15+15+15+15+20"
    27 = "=15+20+15+15+15+20"
    28 = "=14+20+15+15+15+20"
    29 = "=8+15+12+15+12+20"
    30 = "=15+20+15+15+15+20"
    31 = "=15+10+12+12+0+20"
}

foreach ($r in $formulas.Keys) {
    $ws.Range("D$r").Formula = $formulas[$r]
}

# Fix a pre-existing style inconsistency on C31 (it used the even-row
# shaded style even though row 31 is an odd row) so it matches D31/B31.
$ws.Range("C29").Copy()
$ws.Range("C31").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Restore the selected cell
$ws.Range("E4").Select()
